$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43 (pushes existing row 43 and below down by one)
$ws.Rows.Item(43).Insert()

# Fill in the new row 43 with data for LeetCode problem 205 - Isomorphic Strings
$ws.Range("A43").Value = "Hashmap"
$ws.Range("B43").Value = 205
$ws.Range("C43").Value = "205 - Isomorphic Strings"
$ws.Range("D43").Value = "Easy"
$ws.Range("E43").Value = "One pass loop, 2 hashmaps"
$ws.Range("F43").Value = "O(n) time, O(1) memory"
$ws.Range("G43").Value = "O(n) time, O(1) memory"
$ws.Range("H43").Value = "One pass loop, 2 hashmaps"
$ws.Range("I43").Value = "O(n) time, O(1) memory"
$ws.Range("J43").Value = "no"
$ws.Range("K43").Value = "yes"
$ws.Range("L43").Value = "The trick here is really with making sure you check both ways the mapping is the same, and not just one way."
$ws.Range("M43").Value = "35 minutes"
$ws.Range("N43").Value = "YES"

# Match the taller row height used for the new wrapped-text row
$ws.Rows.Item(43).RowHeight = 75

# Extend Table1 to include the newly inserted row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A2:X74"))

# Restore the user's on-screen selection to the newly-added cell
$ws.Range("E43").Select()
